# Fix generation error in the correlation table on Sheet1:
# a batch of coefficients were recomputed / rounded to 3 decimals.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C2"  = 0.935
    "D2"  = 0.908
    "E2"  = 0.944
    "F2"  = 0.942
    "G2"  = 0.939
    "H2"  = 0.945
    "I2"  = 0.909
    "J2"  = 0.943
    "K2"  = 0.941

    "D3"  = 0.915
    "E3"  = 0.94
    "F3"  = 0.945
    "G3"  = 0.934
    "H3"  = 0.942
    "I3"  = 0.918
    "J3"  = 0.944
    "K3"  = 0.949

    "E4"  = 0.917
    "F4"  = 0.91
    "G4"  = 0.926
    "H4"  = 0.913
    "I4"  = 0.95
    "J4"  = 0.917
    "K4"  = 0.92

    "F5"  = 0.94
    "G5"  = 0.938
    "H5"  = 0.945
    "I5"  = 0.914
    "J5"  = 0.938
    "K5"  = 0.942

    "G6"  = 0.934
    "H6"  = 0.943
    "I6"  = 0.907
    "J6"  = 0.942
    "K6"  = 0.941

    "H7"  = 0.938
    "I7"  = 0.92
    "J7"  = 0.939
    "K7"  = 0.941

    "I8"  = 0.909
    "J8"  = 0.946
    "K8"  = 0.947

    "J9"  = 0.917
    "K9"  = 0.924

    "K10" = 0.945
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
